$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '65.462.32'
$ws.Cells.Item(2, 5).Value = '  -4.62%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.258.76'
$ws.Cells.Item(3, 5).Value = '  -5.65%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '''1.00'

# Row 5
$ws.Cells.Item(5, 4).Value = '''554.41'
$ws.Cells.Item(5, 5).Value = '  -3.61%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''179.51'
$ws.Cells.Item(6, 5).Value = '  -5.56%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.05%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.585'
$ws.Cells.Item(8, 5).Value = '  -2.94%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '3.248.90'
$ws.Cells.Item(9, 5).Value = '  -5.55%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -8.37%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''0.584'
$ws.Cells.Item(11, 5).Value = '  -4.70%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''47.27'
$ws.Cells.Item(12, 5).Value = '  -7.30%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''0.0000264'
$ws.Cells.Item(13, 5).Value = '  -6.75%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''629.63'
$ws.Cells.Item(14, 5).Value = '  -1.06%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '''8.53'
$ws.Cells.Item(15, 5).Value = '  -5.58%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '3.766.40'
$ws.Cells.Item(16, 5).Value = '  -6.12%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '65.353.31'
$ws.Cells.Item(17, 5).Value = '  -4.60%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '''17.76'
$ws.Cells.Item(18, 5).Value = '  -1.72%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -3.24%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '3.239.13'
$ws.Cells.Item(20, 5).Value = '  -6.51%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''11.36'
$ws.Cells.Item(21, 5).Value = '  -7.54%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''0.900'
$ws.Cells.Item(22, 5).Value = '  -4.03%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''17.80'
$ws.Cells.Item(23, 5).Value = '  +0.09%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''106.10'
$ws.Cells.Item(24, 5).Value = '  +7.11%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''4.95'
$ws.Cells.Item(25, 5).Value = '  -6.26%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''3.96'
$ws.Cells.Item(26, 5).Value = '  -7.26%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''2.66'
$ws.Cells.Item(27, 5).Value = '  -5.49%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''9.51'
$ws.Cells.Item(28, 5).Value = '  -2.15%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '''8.69'
$ws.Cells.Item(29, 5).Value = '  -5.42%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''30.21'
$ws.Cells.Item(30, 5).Value = '  -6.22%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''3.98'
$ws.Cells.Item(31, 5).Value = '  -3.41%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''6.29'
$ws.Cells.Item(32, 5).Value = '  -6.16%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''11.02'
$ws.Cells.Item(33, 5).Value = '  -4.49%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''546.39'
$ws.Cells.Item(34, 5).Value = '  +9.68%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''0.105'
$ws.Cells.Item(35, 5).Value = '  -2.95%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''56.92'
$ws.Cells.Item(36, 5).Value = '  -6.43%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''1.00'
$ws.Cells.Item(37, 5).Value = '  +0.23%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '3.590.95'
$ws.Cells.Item(38, 5).Value = '  -1.84%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''3.63'
$ws.Cells.Item(39, 5).Value = '  +6.27%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''3.44'
$ws.Cells.Item(40, 5).Value = '  -1.37%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''2.73'
$ws.Cells.Item(41, 5).Value = '  -4.96%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''0.129'
$ws.Cells.Item(42, 5).Value = '  -1.72%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '0.0₃0707'
$ws.Cells.Item(43, 5).Value = '  -8.25%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''31.98'
$ws.Cells.Item(44, 5).Value = '  -6.68%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '''0.336'
$ws.Cells.Item(45, 5).Value = '  -8.14%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''3.29'
$ws.Cells.Item(46, 5).Value = '  -1.44%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''0.0412'
$ws.Cells.Item(47, 5).Value = '  -5.33%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'ThetaToken'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(48, 4).Value = '''2.60'
$ws.Cells.Item(48, 5).Value = '  -6.70%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Stellar'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(49, 4).Value = '''0.129'
$ws.Cells.Item(49, 5).Value = '  -3.43%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -0.23%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''1.25'
$ws.Cells.Item(51, 5).Value = '  +1.96%  '
